$wb = $excel.ActiveWorkbook

# --- "esneklik" sheet: update elasticity values and add a new TCO row ---
$ws = $wb.Worksheets.Item("esneklik")

# B2 (Kendi fiyat esnekligi) changes from -1.66 to -2.8
$ws.Range("B2").Value = -2.8

# B3 (Rakip fiyat esnekligi) changes from 0.82 to 0
$ws.Range("B3").Value = 0

# New row 4: TCO-talep esnekligi / -2.8
$ws.Range("A4").Value = "TCO-talep esnekligi"
$ws.Range("B4").Value = -2.8

# Make "esneklik" the active sheet/tab and select B6 on it
$ws.Activate()
$ws.Range("B6").Select()
